$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B4").Value = "SingleUseId1"
$ws.Range("C4").Value = "Default"
$ws.Range("D4").Value = "Center"
$ws.Range("E4").Value = "LTR"
$ws.Range("F4").Value = "Toggle LED"
